# The deck's theme colour scheme (the clrScheme living in the theme part
# bound to the slide master / slides) is being swapped from the custom
# "Integral" / "Red Violet" palette back to the stock PowerPoint
# "Office Theme" palette.
#
# PowerPoint's legacy 8-slot ColorScheme object (Slide.ColorScheme /
# NotesMaster.ColorScheme / HandoutMaster.ColorScheme - they all resolve to
# the presentation's single master colour scheme) exposes 12 indexed RGB
# slots that line up 1:1 with the DrawingML <a:clrScheme> children:
#   1 dk1   2 lt1   3 dk2   4 lt2
#   5 accent1  6 accent2  7 accent3  8 accent4
#   9 accent5  10 accent6  11 hlink  12 folHlink
#
# RGBColor.RGB takes a COM OLE_COLOR style integer (R + G*256 + B*65536),
# i.e. the same encoding VBA's RGB(r,g,b) produces, so a small helper turns
# a plain "RRGGBB" hex string into the right integer.

function HexToRgbVal($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$cs = $s.ColorScheme

$cs.Colors(1).RGB  = HexToRgbVal "000000"   # dk1
$cs.Colors(2).RGB  = HexToRgbVal "FFFFFF"   # lt1
$cs.Colors(3).RGB  = HexToRgbVal "44546A"   # dk2
$cs.Colors(4).RGB  = HexToRgbVal "E7E6E6"   # lt2
$cs.Colors(5).RGB  = HexToRgbVal "5B9BD5"   # accent1
$cs.Colors(6).RGB  = HexToRgbVal "ED7D31"   # accent2
$cs.Colors(7).RGB  = HexToRgbVal "A5A5A5"   # accent3
$cs.Colors(8).RGB  = HexToRgbVal "FFC000"   # accent4
$cs.Colors(9).RGB  = HexToRgbVal "4472C4"   # accent5
$cs.Colors(10).RGB = HexToRgbVal "70AD47"   # accent6
$cs.Colors(11).RGB = HexToRgbVal "0563C1"   # hlink
$cs.Colors(12).RGB = HexToRgbVal "954F72"   # folHlink
